# Update the "rand_digit" column (J) values on Sheet1 to match the new data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "J2"  = 5
    "J3"  = 2
    "J4"  = 7
    "J5"  = 2
    "J6"  = 1
    "J7"  = 5
    "J8"  = 1
    "J9"  = 2
    "J10" = 8
    "J11" = 8
    "J12" = 6
    "J14" = 7
    "J15" = 4
    "J16" = 1
    "J17" = 3
    "J18" = 1
    "J19" = 1
    "J21" = 1
    "J22" = 4
    "J23" = 7
    "J24" = 1
    "J25" = 8
    "J26" = 5
    "J27" = 6
    "J28" = 6
    "J29" = 2
    "J31" = 6
    "J32" = 7
    "J33" = 8
    "J34" = 2
    "J35" = 5
    "J36" = 2
    "J38" = 4
    "J39" = 4
    "J40" = 6
    "J41" = 1
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
